$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.591.28'
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("D3").Value = '2.286.25'
$ws.Range("E3").Value = '  +4.63%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.12'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.56%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.635'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0969'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.61%  '
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").Value = '2.622.06'
$ws.Range("E15").Value = '  +4.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.881'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.17%  '
$ws.Range("D18").Value = '2.272.03'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '42.564.96'
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("D20").Value = '0.0₃0997'
$ws.Range("E20").Value = '  +5.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.42%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0806'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.17%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +23.04%  '
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.126'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.51%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.75'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0309'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +19.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.213'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.104'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.62%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.02%  '
